# Mark "UploadCountryVersion" (row 9) as DONE instead of TODO,
# matching the already-completed rows (e.g. B37:B39 -> green "DONE"),
# and move the active selection from D24 to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status text for UploadCountryVersion (B9): "TODO" -> "DONE"
$ws.Range("B9").Value = "DONE"

# Re-color the cell to match the green "DONE" fill used elsewhere (e.g. B37:B39)
$ws.Range("B9").Interior.Color = 5287936

# Move the current selection to D11 (was D24)
$ws.Range("D11").Select()
